$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CL")

# Row 4 - Inventory
$ws.Range("C4").Value = 1673000000.0
$ws.Range("D4").Value = 1578000000.0
$ws.Range("E4").Value = 1524000000.0
$ws.Range("F4").Value = 1301000000.0
$ws.Range("G4").Value = 1400000000.0

# Row 10 - Long Term Assets (Tax, Deferred)
$ws.Range("B10").Value = 201000000.0

# Row 11 - Long-term assets (Other)
$ws.Range("B11").Value = 940000000.0

# Row 14 - Notes Payable
$ws.Range("B14").Value = 254000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 1306000000.0
$ws.Range("C15").Value = 1393000000.0
$ws.Range("D15").Value = 1225000000.0
$ws.Range("E15").Value = 1189000000.0
$ws.Range("F15").Value = 1216000000.0
$ws.Range("G15").Value = 1237000000.0

# Row 16 - Current Part of Debt
$ws.Range("B16").Value = 9000000.0

# Row 17 - Current Part of Taxes to Pay
$ws.Range("B17").Value = 422000000.0

# Row 19 - Total current liabilities
$ws.Range("B19").Value = 4539000000.0

# Row 20 - Long Term Debt (Total)
$ws.Range("B20").Value = 7570000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 415000000.0
$ws.Range("C21").Value = 135000000.0
$ws.Range("D21").Value = 198000000.0
$ws.Range("E21").Value = 176000000.0
$ws.Range("F21").Value = 236000000.0
$ws.Range("G21").Value = 330000000.0

# Row 22 - Non-current Liabilities (Other)
$ws.Range("B22").Value = 2614000000.0

# Row 23 - Total non-current liabilities
$ws.Range("B23").Value = 10599000000.0

# Row 25 - Additional Paid In Capital
$ws.Range("B25").Value = 3011000000.0

# Row 26 - Common Stock (Net)
$ws.Range("B26").Value = 1466000000.0

# Row 27 - Retained Earnings
$ws.Range("B27").Value = 23624000000.0

# Row 28 - Treasury Stock
$ws.Range("B28").Value = 23384000000.0

# Row 29 - Common Equity (Total)
$ws.Range("B29").Value = 663000000.0

# Row 32 - Shares (Common)
$ws.Range("B32").Value = 845969000.0

# Row 34 - Net Debt
$ws.Range("B34").Value = 6838000000.0

# Row 35 - Total Debt
$ws.Range("B35").Value = 7833000000.0
